$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.26126961572082763
$ws.Range("A2").Value = -0.0059999999251196812
$ws.Range("A3").Value = -0.0039999999319597634
$ws.Range("A4").Value = -0.0079999998776809633
$ws.Range("A5").Value = -0.002999999929189201
$ws.Range("A6").Value = -0.0019999999242834576
$ws.Range("A7").Value = -0.0099999998300019932
$ws.Range("A8").Value = -0.0099999998271713686
$ws.Range("A9").Value = -0.0019999999200894791
$ws.Range("A10").Value = -0.0019999999190591922
$ws.Range("A11").Value = -0.0029999999073444528
$ws.Range("A12").Value = -0.0034999999014897476
$ws.Range("A13").Value = -0.0034999999022389261
$ws.Range("A14").Value = -0.0079999998499395986
$ws.Range("A15").Value = -0.00099999993381860719
$ws.Range("A16").Value = -0.0019999999239397326
$ws.Range("A17").Value = 0.017719031453101053
$ws.Range("A18").Value = -0.0047988707352688564
$ws.Range("A19").Value = -0.039171490531437758
$ws.Range("A20").Value = -0.0039999999270818876
$ws.Range("A21").Value = -0.044705075393508764
$ws.Range("A22").Value = -0.0039999999225273086
$ws.Range("A23").Value = -0.0049999999159098252
$ws.Range("A24").Value = -0.019999999732645435
$ws.Range("A25").Value = -0.019999999728833373
$ws.Range("A26").Value = -0.0024999999183918931
$ws.Range("A27").Value = -0.0024999999150385754
$ws.Range("A28").Value = 0.032173019321861673
$ws.Range("A29").Value = -0.006999999833633197
$ws.Range("A30").Value = -0.059999999200295306
$ws.Range("A31").Value = -0.0069999998206142777
$ws.Range("A32").Value = -0.0099999997848172484
$ws.Range("A33").Value = -0.003999999855475167
